$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.21%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.99%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.128"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.10%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08078"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.44%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.442"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.78%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.943"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.50%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.296"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.01%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.00%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9382"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.03%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1319"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.93%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1937"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.15%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09090"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.18%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03484"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.32%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09676"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.46%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001411"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.39%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006142"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.76%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.573"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-5.65%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.39%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1282"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.32%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.024"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "7.78%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04357"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.43%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001243"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004724"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003786"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "190.72%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02213"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.49%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05243"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.74%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007590"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.00%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01032"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.65%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1387"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.87%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002036"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.56%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009098"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.49%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006605"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.85%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003012"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "18.01%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.01%"
